$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "South Korea"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "4"
$ws.Range("C2").Value = "Insurance (Life)"

$ws.Range("A3").Value = "South Korea"
$ws.Range("B3").Value = "Tong Yang Life Insurance Co., Ltd. (KOSE:A082640)"
$ws.Range("C3").Value = "Insurance (Life)"

$ws.Range("A4").Value = "South Korea"
$ws.Range("B4").Value = "Mirae Asset Life Insurance Co., Ltd. (KOSE:A085620)"
$ws.Range("C4").Value = "Insurance (Life)"

$ws.Range("A5").Value = "South Korea"
$ws.Range("B5").Value = "Hanwha Life Insurance Co., Ltd. (KOSE:A088350)"
$ws.Range("C5").Value = "Insurance (Life)"

$ws.Range("A6").Value = "South Korea"
$ws.Range("B6").Value = "Samsung Life Insurance Co., Ltd. (KOSE:A032830)"
$ws.Range("C6").Value = "Insurance (Life)"

$row = New-Object 'object[,]' 1,40
$row[0,0] = 0.05385
$row[0,1] = -0.06155
$row[0,2] = 0.006
$row[0,3] = 0.03551234325122792
$row[0,4] = 0.03551234325122792
$row[0,5] = 0.0585286881122463
$row[0,6] = 0.04422550048259027
$row[0,7] = 1255.5
$row[0,8] = 0.02296695356303336
$row[0,9] = 495.2395
$row[0,10] = 0.03151340740175117
$row[0,11] = 0.3944559936280366
$row[0,12] = 479.5394999999999
$row[0,13] = 0.03051437461820403
$row[0,14] = 0.3819510155316606
$row[0,15] = 15.7
$row[0,16] = 0.03170183315345405
$row[0,17] = 1.82
$row[0,18] = 0.0001158114436978212
$row[0,19] = 0.03877158761062589
$row[0,20] = 0.07599041620467631
$row[0,21] = -0.03721882859405042
$row[0,22] = 0.9371957870413025
$row[0,23] = 0.07255998954774698
$row[0,24] = 0.05969261982967428
$row[0,25] = 0.01355557142778391
$row[0,26] = 17123.6
$row[0,27] = 0
$row[0,28] = 17123.6
$row[0,29] = 17121.78
$row[0,30] = 0.5214441453402683
$row[0,31] = 0.2502429564548964
$row[0,32] = 0.5214176212306979
$row[0,33] = 0.2502230143828201
$row[0,34] = 483.54
$row[0,35] = 483.54
$row[0,36] = 4.330374528993754
$row[0,37] = 6.616825908921702
$row[0,38] = 4.329914270540931
$row[0,39] = 6.616825908921702
$ws.Range("D2:AQ2").Value = $row

$row = New-Object 'object[,]' 1,40
$row[0,0] = 0.0697
$row[0,1] = -0.0687
$row[0,2] = 0.006
$row[0,3] = 0.07677954662749613
$row[0,4] = 0.07677954662749613
$row[0,5] = 0.0729485092699194
$row[0,6] = 0.05134175758207304
$row[0,7] = 99.8
$row[0,8] = 0.01981023462622573
$row[0,9] = 31.0241
$row[0,10] = 0.06153133677112258
$row[0,11] = 0.3108627254509019
$row[0,12] = 31.0241
$row[0,13] = 0.06153133677112258
$row[0,14] = 0.3108627254509019
$row[0,15] = 0
$row[0,16] = 0
$row[0,17] = 0
$row[0,18] = 0
$row[0,19] = 0.04901768172888016
$row[0,20] = 0.07886467302683596
$row[0,21] = -0.0298469912979558
$row[0,22] = 2.369839467343557
$row[0,23] = 0.121671723440782
$row[0,24] = 0.05974003113166211
$row[0,25] = 0.06193169230911989
$row[0,26] = 259.9
$row[0,27] = 0
$row[0,28] = 259.9
$row[0,29] = 259.9
$row[0,30] = 0.3401387252977359
$row[0,31] = 0.1034222045364107
$row[0,32] = 0.3401387252977359
$row[0,33] = 0.1034222045364107
$row[0,34] = 12
$row[0,35] = 12
$row[0,36] = 0.6719234746639089
$row[0,37] = 30.625
$row[0,38] = 0.6719234746639089
$row[0,39] = 30.625
$ws.Range("D3:AQ3").Value = $row

$ws.Range("F4").ClearContents()

$row = New-Object 'object[,]' 1,2
$row[0,0] = -0.0352
$row[0,1] = 0.05980000000000001
$ws.Range("D4:E4").Value = $row
$row = New-Object 'object[,]' 1,37
$row[0,0] = 0.04495806150978564
$row[0,1] = 0.04495806150978564
$row[0,2] = 0.08354147250698975
$row[0,3] = 0.06450983527906605
$row[0,4] = 104.4
$row[0,5] = 0.03891891891891892
$row[0,6] = 35.191
$row[0,7] = 0.07450984543722211
$row[0,8] = 0.3370785440613027
$row[0,9] = 19.491
$row[0,10] = 0.04126826169807326
$row[0,11] = 0.1866954022988506
$row[0,12] = 15.7
$row[0,13] = 0.4461367963399733
$row[0,14] = 0
$row[0,15] = 0
$row[0,16] = 0.05526148634342579
$row[0,17] = 0.07311615938251664
$row[0,18] = -0.01785467303909085
$row[0,19] = 1.251401849975462
$row[0,20] = 0.08072772720983555
$row[0,21] = 0.05836362771223969
$row[0,22] = 0.02236409949759586
$row[0,23] = 171.3
$row[0,24] = 0
$row[0,25] = 171.3
$row[0,26] = 171.3
$row[0,27] = 0.2661591050341827
$row[0,28] = 0.07601508764144664
$row[0,29] = 0.2661591050341827
$row[0,30] = 0.07601508764144664
$row[0,31] = 8.44
$row[0,32] = 8.44
$row[0,33] = 0.6290855673889093
$row[0,34] = 26.5521327014218
$row[0,35] = 0.6290855673889093
$row[0,36] = 26.5521327014218
$ws.Range("G4:AQ4").Value = $row

$ws.Range("AO5").ClearContents()
$ws.Range("AQ5").ClearContents()

$row = New-Object 'object[,]' 1,37
$row[0,0] = 0.0925
$row[0,1] = -0.162
$row[0,2] = 0.005
$row[0,3] = -0.01313980368357701
$row[0,4] = -0.01313980368357701
$row[0,5] = 0.04288944736272204
$row[0,6] = 0.03348927654640431
$row[0,7] = 196.1
$row[0,8] = 0.009876257195666736
$row[0,9] = 19.5364
$row[0,10] = 0.01159499080064099
$row[0,11] = 0.09962468128505864
$row[0,12] = 19.5364
$row[0,13] = 0.01159499080064099
$row[0,14] = 0.09962468128505864
$row[0,15] = 0
$row[0,16] = 0
$row[0,17] = 0
$row[0,18] = 0
$row[0,19] = 0.01988642125545077
$row[0,20] = 0.07014822234754778
$row[0,21] = -0.05026180109209701
$row[0,22] = 1.922772258052021
$row[0,23] = 0.06439225188565841
$row[0,24] = 0.05964520852768644
$row[0,25] = 0.004747043357971964
$row[0,26] = 478.2
$row[0,27] = 0
$row[0,28] = 478.2
$row[0,29] = 478.2
$row[0,30] = 0.2210716101890805
$row[0,31] = 0.03521717997437144
$row[0,32] = 0.2210716101890805
$row[0,33] = 0.03521717997437144
$row[0,34] = 0
$row[0,35] = 0
$row[0,36] = 0.4355984696666059
$ws.Range("D5:AN5").Value = $row
$ws.Range("AP5").Value = 0.4355984696666059

$row = New-Object 'object[,]' 1,40
$row[0,0] = 0.038
$row[0,1] = -0.0544
$row[0,2] = 0.105
$row[0,3] = 0.06256298565865002
$row[0,4] = 0.06256298565865002
$row[0,5] = 0.06483323797043135
$row[0,6] = 0.04964010141906125
$row[0,7] = 855.2
$row[0,8] = 0.03156942726886802
$row[0,9] = 409.4879999999999
$row[0,10] = 0.03136925646171996
$row[0,11] = 0.4788213283442468
$row[0,12] = 409.4879999999999
$row[0,13] = 0.03136925646171996
$row[0,14] = 0.4788213283442468
$row[0,15] = 0
$row[0,16] = 0
$row[0,17] = 1.82
$row[0,18] = 0.0001394230032634176
$row[0,19] = 0.02852549349237163
$row[0,20] = 0.106205695319566
$row[0,21] = -0.07768020182719439
$row[0,22] = 0.6194320967328869
$row[0,23] = 0.03074867210404227
$row[0,24] = 0.05991033845548902
$row[0,25] = -0.02916166635144675
$row[0,26] = 16214.2
$row[0,27] = 0
$row[0,28] = 16214.2
$row[0,29] = 16212.38
$row[0,30] = 0.5539907065737324
$row[0,31] = 0.3237478735214485
$row[0,32] = 0.5539629702270676
$row[0,33] = 0.3237232977469691
$row[0,34] = 463.1
$row[0,35] = 463.1
$row[0,36] = 7.378811322472012
$row[0,37] = 3.792485424314403
$row[0,38] = 7.377983070901975
$row[0,39] = 3.792485424314403
$ws.Range("D6:AQ6").Value = $row

$ws.Rows.Item(7).Delete()

Write-Host "Edit complete"